# Commit: fix: modify description
# 1) Update the cached "datetimeFigureOut" date field text from 2019/1/4 to
#    2019/1/20 everywhere it appears: the slide master, every slide layout,
#    and the notes master.
# 2) Remove the stray "Dijkstra" run that preceded "最短路径算法" in the
#    bullet list on slide 11.

$p = $ppt.ActivePresentation

$oldDate = "2019/1/4"
$newDate = "2019/1/20"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        if (-not $shp.TextFrame.HasText) { continue }

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if (-not $isDatePlaceholder) { continue }

        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide 11: drop the leading "Dijkstra" run from the "最短路径算法" bullet.
$s11 = $p.Slides.Item(11)
$shp = $s11.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange
$c = $tr.Characters(1, 8)
if ($c.Text -eq "Dijkstra") {
    $c.Text = ""
}
